# "try to fix the spider attack after dead, but not succes"
#
# - Adds a new Todo row (#24): "All UI item should be private" / Fish / 7jan
# - Adds a new Defect row (#6): Enemy attack / Enemy could still attack once
#   after its dead... / Forest / Fish / 7jan
# - Re-labels every existing Create/Finish date cell (column F/G, both
#   sheets) from a numeric date to the literal text "7jan", and switches
#   those columns to a text number format
# - Widens column A on the Todo sheet, and widens the F/G date columns
# - Moves the active sheet/selection from Todo!F24 to Defect!B15 (selecting
#   Todo!G12 along the way)
# - Sets a portrait page setup on the Todo sheet

$wb = $excel.ActiveWorkbook

$todo = $wb.Worksheets.Item("Todo ")
$defect = $wb.Worksheets.Item("Defect")

# ---------------------------------------------------------------------
# 1. Re-label every existing Create date / Finish date cell as "7jan"
#    and give it a text number format (column F/G, only cells that
#    actually hold a date already).
# ---------------------------------------------------------------------

$todoDateCells = @("F2","F3","F4","G4","F5","F6","F7","F8","F9","G9","F10","F11","G11","F12","F13","F14","F15","F16","F17","F18","F19","F20","F21","F22","F23")
foreach ($addr in $todoDateCells) {
    $todo.Range($addr).Value = "7jan"
    $todo.Range($addr).NumberFormat = "@"
}

$defectDateCells = @("F2","G2","F3","F4","F5")
foreach ($addr in $defectDateCells) {
    $defect.Range($addr).Value = "7jan"
    $defect.Range($addr).NumberFormat = "@"
}

# The "Create date" / "Finish date" header cells pick up the same text
# format as the column below them.
$todo.Range("F1:G1").NumberFormat = "@"
$defect.Range("F1:G1").NumberFormat = "@"

# ---------------------------------------------------------------------
# 2. New rows
# ---------------------------------------------------------------------

# Todo sheet gained a new issue: "All UI item should be private"
$todo.Range("A24").Value = "All UI item should be private"
$todo.Range("D24").Value = "Fish"
$todo.Range("F24").Value = "7jan"
$todo.Range("F24").NumberFormat = "@"

# Defect sheet gained a new defect: "Enemy attack"
$defect.Range("A6").Value = "Enemy attack"
$defect.Range("B6").Value = "Enemy could still attack once after its dead. It leads to the level up times cut to once only"
$defect.Range("C6").Value = "Forest"
$defect.Range("D6").Value = "Fish"
$defect.Range("F6").Value = "7jan"
$defect.Range("F6").NumberFormat = "@"

# ---------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------

# (ColumnWidth is stored by the engine quantised to 1/7 character widths plus
# a fixed 5/7 padding, so these inputs are chosen to land as close as
# possible to the target stored widths of 25.296875 / 13.59765625.)
$todo.Columns.Item(1).ColumnWidth = 24.571428571428573
$todo.Columns.Item(6).ColumnWidth = 12.857142857142858

# ---------------------------------------------------------------------
# 4. Page setup
# ---------------------------------------------------------------------

$todo.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5. Selection / active sheet: end up on Defect!B15, having passed
#    through Todo!G12.
# ---------------------------------------------------------------------

$todo.Range("G12").Select()
$defect.Activate()
$defect.Range("B15").Select()
